$wb = $excel.ActiveWorkbook

# ---- Sheet 1: appending rows 193-203 ----
$ws = $wb.Worksheets.Item(1)
$startRow = 193
$endRow = 203
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"
$ws.Range("C$startRow`:C$endRow").NumberFormat = "@"
$ws.Range("I$startRow`:I$endRow").NumberFormat = "@"
$data = New-Object 'object[,]' 11,9
$data[0,0] = '大智 (稳健智远)'
$data[0,1] = '000333'
$data[0,2] = '美的集团'
$data[0,3] = 3.06
$data[0,4] = 43.36476491361748
$data[0,5] = 72.2
$data[0,6] = 3130.936026763182
$data[0,7] = 102355.8645050981
$data[0,8] = '202506301630'
$data[1,0] = '大智 (稳健智远)'
$data[1,1] = '510050'
$data[1,2] = '上证50ETF'
$data[1,3] = 5.14
$data[1,4] = 1872.352431685969
$data[1,5] = 2.81
$data[1,6] = 5261.310333037572
$data[1,7] = 102355.8645050981
$data[1,8] = '202506301630'
$data[2,0] = '大智 (稳健智远)'
$data[2,1] = '510300'
$data[2,2] = '沪深300ETF'
$data[2,3] = 5.01
$data[2,4] = 1287.543601270288
$data[2,5] = 3.98
$data[2,6] = 5124.423533055747
$data[2,7] = 102355.8645050981
$data[2,8] = '202506301630'
$data[3,0] = '大智 (稳健智远)'
$data[3,1] = '518880'
$data[3,2] = '黄金ETF'
$data[3,3] = 4.87
$data[3,4] = 681.5657355049799
$data[3,5] = 7.31
$data[3,6] = 4982.245526541403
$data[3,7] = 102355.8645050981
$data[3,8] = '202506301630'
$data[4,0] = '大智 (稳健智远)'
$data[4,1] = '600085'
$data[4,2] = '同仁堂'
$data[4,3] = 1.98
$data[4,4] = 56.17257911411885
$data[4,5] = 36.06
$data[4,6] = 2025.583202855126
$data[4,7] = 102355.8645050981
$data[4,8] = '202506301630'
$data[5,0] = '大智 (稳健智远)'
$data[5,1] = '600900'
$data[5,2] = '长江电力'
$data[5,3] = 30.04
$data[5,4] = 1020.137146333154
$data[5,5] = 30.14
$data[5,6] = 30746.93359048126
$data[5,7] = 102355.8645050981
$data[5,8] = '202506301630'
$data[6,0] = '大智 (稳健智远)'
$data[6,1] = '600989'
$data[6,2] = '宝丰能源'
$data[6,3] = 4.8
$data[6,4] = 304.275626545359
$data[6,5] = 16.14
$data[6,6] = 4911.008612442094
$data[6,7] = 102355.8645050981
$data[6,8] = '202506301630'
$data[7,0] = '大智 (稳健智远)'
$data[7,1] = '601899'
$data[7,2] = '紫金矿业'
$data[7,3] = 9.91
$data[7,4] = 520.2217600143598
$data[7,5] = 19.5
$data[7,6] = 10144.32432028002
$data[7,7] = 102355.8645050981
$data[7,8] = '202506301630'
$data[8,0] = '大智 (稳健智远)'
$data[8,1] = 'HK02899'
$data[8,2] = '紫金矿业'
$data[8,3] = 9.890000000000001
$data[8,4] = 504.6340243688386
$data[8,5] = 20.05
$data[8,6] = 10117.91218859521
$data[8,7] = 102355.8645050981
$data[8,8] = '202506301630'
$data[9,0] = '大智 (稳健智远)'
$data[9,1] = 'HK06881'
$data[9,2] = '中国银河'
$data[9,3] = 5.22
$data[9,4] = 604.9539788562929
$data[9,5] = 8.83
$data[9,6] = 5341.743633301066
$data[9,7] = 102355.8645050981
$data[9,8] = '202506301630'
$data[10,0] = '大智 (稳健智远)'
$data[10,1] = '100000'
$data[10,2] = '现金'
$data[10,3] = 20.1
$data[10,4] = 20569.44353774537
$data[10,5] = 1
$data[10,6] = 20569.44353774537
$data[10,7] = 102355.8645050981
$data[10,8] = '202506301630'
$ws.Range("A$startRow`:I$endRow").Value = $data

# ---- Sheet 2: appending rows 130-137 ----
$ws = $wb.Worksheets.Item(2)
$startRow = 130
$endRow = 137
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"
$ws.Range("C$startRow`:C$endRow").NumberFormat = "@"
$ws.Range("I$startRow`:I$endRow").NumberFormat = "@"
$data = New-Object 'object[,]' 8,9
$data[0,0] = '大成 (锐进先锋)'
$data[0,1] = '000725'
$data[0,2] = '京东方A'
$data[0,3] = 4.84
$data[0,4] = 1248.221835380318
$data[0,5] = 3.99
$data[0,6] = 4980.405123167469
$data[0,7] = 102805.3771514609
$data[0,8] = '202506301630'
$data[1,0] = '大成 (锐进先锋)'
$data[1,1] = '159781'
$data[1,2] = '科创创业ETF'
$data[1,3] = 9.960000000000001
$data[1,4] = 18275.64332089199
$data[1,5] = 0.5600000000000001
$data[1,6] = 10234.36025969952
$data[1,7] = 102805.3771514609
$data[1,8] = '202506301630'
$data[2,0] = '大成 (锐进先锋)'
$data[2,1] = '513100'
$data[2,2] = '纳指ETF'
$data[2,3] = 4.87
$data[2,4] = 3071.323391427681
$data[2,5] = 1.63
$data[2,6] = 5006.25712802712
$data[2,7] = 102805.3771514609
$data[2,8] = '202506301630'
$data[3,0] = '大成 (锐进先锋)'
$data[3,1] = '513290'
$data[3,2] = '纳指生物科技ETF'
$data[3,3] = 0.93
$data[3,4] = 860.2723995645101
$data[3,5] = 1.11
$data[3,6] = 954.9023635166063
$data[3,7] = 102805.3771514609
$data[3,8] = '202506301630'
$data[4,0] = '大成 (锐进先锋)'
$data[4,1] = '603119'
$data[4,2] = '浙江荣泰'
$data[4,3] = 42.8
$data[4,4] = 951.4702743128352
$data[4,5] = 46.24
$data[4,6] = 43995.9854842255
$data[4,7] = 102805.3771514609
$data[4,8] = '202506301630'
$data[5,0] = '大成 (锐进先锋)'
$data[5,1] = '688290'
$data[5,2] = '景业智能'
$data[5,3] = 7.46
$data[5,4] = 147.1698966627607
$data[5,5] = 52.12
$data[5,6] = 7670.495014063087
$data[5,7] = 102805.3771514609
$data[5,8] = '202506301630'
$data[6,0] = '大成 (锐进先锋)'
$data[6,1] = 'HK01896'
$data[6,2] = '猫眼娱乐'
$data[6,3] = 0.97
$data[6,4] = 136.3853979171044
$data[6,5] = 7.29
$data[6,6] = 994.2495508156911
$data[6,7] = 102805.3771514609
$data[6,8] = '202506301630'
$data[7,0] = '大成 (锐进先锋)'
$data[7,1] = '100000'
$data[7,2] = '现金'
$data[7,3] = 28.18
$data[7,4] = 28968.72222794588
$data[7,5] = 1
$data[7,6] = 28968.72222794588
$data[7,7] = 102805.3771514609
$data[7,8] = '202506301630'
$ws.Range("A$startRow`:I$endRow").Value = $data

# ---- Sheet 3: appending rows 273-288 ----
$ws = $wb.Worksheets.Item(3)
$startRow = 273
$endRow = 288
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"
$ws.Range("C$startRow`:C$endRow").NumberFormat = "@"
$ws.Range("I$startRow`:I$endRow").NumberFormat = "@"
$data = New-Object 'object[,]' 16,9
$data[0,0] = '范式进化投资组合'
$data[0,1] = '000333'
$data[0,2] = '美的集团'
$data[0,3] = 1.02
$data[0,4] = 14.22471587925232
$data[0,5] = 72.2
$data[0,6] = 1027.024486482018
$data[0,7] = 100822.2970584788
$data[0,8] = '202506301630'
$data[1,0] = '范式进化投资组合'
$data[1,1] = '000725'
$data[1,2] = '京东方A'
$data[1,3] = 5.08
$data[1,4] = 1282.574489616538
$data[1,5] = 3.99
$data[1,6] = 5117.472213569987
$data[1,7] = 100822.2970584788
$data[1,8] = '202506301630'
$data[2,0] = '范式进化投资组合'
$data[2,1] = '159781'
$data[2,2] = '科创创业ETF'
$data[2,3] = 5.25
$data[2,4] = 9459.512127702146
$data[2,5] = 0.5600000000000001
$data[2,6] = 5297.326791513203
$data[2,7] = 100822.2970584788
$data[2,8] = '202506301630'
$data[3,0] = '范式进化投资组合'
$data[3,1] = '510050'
$data[3,2] = '上证50ETF'
$data[3,3] = 5.14
$data[3,4] = 1843.714434635192
$data[3,5] = 2.81
$data[3,6] = 5180.83756132489
$data[3,7] = 100822.2970584788
$data[3,8] = '202506301630'
$data[4,0] = '范式进化投资组合'
$data[4,1] = '510300'
$data[4,2] = '沪深300ETF'
$data[4,3] = 5
$data[4,4] = 1265.443994258174
$data[4,5] = 3.98
$data[4,6] = 5036.467097147533
$data[4,7] = 100822.2970584788
$data[4,8] = '202506301630'
$data[5,0] = '范式进化投资组合'
$data[5,1] = '513100'
$data[5,2] = '纳指ETF'
$data[5,3] = 1.02
$data[5,4] = 628.619842441686
$data[5,5] = 1.63
$data[5,6] = 1024.650343179948
$data[5,7] = 100822.2970584788
$data[5,8] = '202506301630'
$data[6,0] = '范式进化投资组合'
$data[6,1] = '513290'
$data[6,2] = '纳指生物科技ETF'
$data[6,3] = 0.98
$data[6,4] = 890.1943286597277
$data[6,5] = 1.11
$data[6,6] = 988.1157048122978
$data[6,7] = 100822.2970584788
$data[6,8] = '202506301630'
$data[7,0] = '范式进化投资组合'
$data[7,1] = '518880'
$data[7,2] = '黄金ETF'
$data[7,3] = 0.98
$data[7,4] = 134.8043253495631
$data[7,5] = 7.31
$data[7,6] = 985.4196183053062
$data[7,7] = 100822.2970584788
$data[7,8] = '202506301630'
$data[8,0] = '范式进化投资组合'
$data[8,1] = '600085'
$data[8,2] = '同仁堂'
$data[8,3] = 0.99
$data[8,4] = 27.5486702579541
$data[8,5] = 36.06
$data[8,6] = 993.4050495018249
$data[8,7] = 100822.2970584788
$data[8,8] = '202506301630'
$data[9,0] = '范式进化投资组合'
$data[9,1] = '600900'
$data[9,2] = '长江电力'
$data[9,3] = 9.970000000000001
$data[9,4] = 333.5366627327556
$data[9,5] = 30.14
$data[9,6] = 10052.79501476525
$data[9,7] = 100822.2970584788
$data[9,8] = '202506301630'
$data[10,0] = '范式进化投资组合'
$data[10,1] = '600989'
$data[10,2] = '宝丰能源'
$data[10,3] = 4.8
$data[10,4] = 299.6922357452684
$data[10,5] = 16.14
$data[10,6] = 4837.032684928632
$data[10,7] = 100822.2970584788
$data[10,8] = '202506301630'
$data[11,0] = '范式进化投资组合'
$data[11,1] = '601899'
$data[11,2] = '紫金矿业'
$data[11,3] = 9.9
$data[11,4] = 511.7946231960747
$data[11,5] = 19.5
$data[11,6] = 9979.995152323458
$data[11,7] = 100822.2970584788
$data[11,8] = '202506301630'
$data[12,0] = '范式进化投资组合'
$data[12,1] = '603119'
$data[12,2] = '浙江荣泰'
$data[12,3] = 1.13
$data[12,4] = 24.54234478152828
$data[12,5] = 46.24
$data[12,6] = 1134.838022697868
$data[12,7] = 100822.2970584788
$data[12,8] = '202506301630'
$data[13,0] = '范式进化投资组合'
$data[13,1] = 'HK01896'
$data[13,2] = '猫眼娱乐'
$data[13,3] = 0.2
$data[13,4] = 27.36177290206568
$data[13,5] = 7.29
$data[13,6] = 199.4673244560588
$data[13,7] = 100822.2970584788
$data[13,8] = '202506301630'
$data[14,0] = '范式进化投资组合'
$data[14,1] = 'HK06881'
$data[14,2] = '中国银河'
$data[14,3] = 1.05
$data[14,4] = 119.3568656159162
$data[14,5] = 8.83
$data[14,6] = 1053.92112338854
$data[14,7] = 100822.2970584788
$data[14,8] = '202506301630'
$data[15,0] = '范式进化投资组合'
$data[15,1] = '100000'
$data[15,2] = '现金'
$data[15,3] = 47.52
$data[15,4] = 47913.52887008203
$data[15,5] = 1
$data[15,6] = 47913.52887008203
$data[15,7] = 100822.2970584788
$data[15,8] = '202506301630'
$ws.Range("A$startRow`:I$endRow").Value = $data

Write-Host "Appended rows to all three sheets successfully."